$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 19931.8
$ws.Range("J51").Value = 18636.182
$ws.Range("L51").Value = 18636.182
$ws.Range("N51").Value = -19604.182

$ws.Range("H53").Value = 464
$ws.Range("I53").Value = 330.9091
$ws.Range("J53").Value = 626.6667
$ws.Range("K53").Value = 330.9091
$ws.Range("L53").Value = 626.6667
$ws.Range("M53").Value = 306.0909
$ws.Range("N53").Value = -1900.6667

$ws.Range("H98").Value = 2019978.4
$ws.Range("I98").Value = 2116036.5
$ws.Range("J98").Value = 1429906.6
$ws.Range("K98").Value = 2116036.5
$ws.Range("L98").Value = 1429906.6
$ws.Range("M98").Value = -2114538.5
$ws.Range("N98").Value = -1432902.6

$ws.Range("H100").Value = 1531.6666
$ws.Range("I100").Value = 1471.3334
$ws.Range("J100").Value = 1833.3334
$ws.Range("K100").Value = 1471.3334
$ws.Range("L100").Value = 1833.3334
$ws.Range("M100").Value = -930.3334
$ws.Range("N100").Value = -2915.3334

$ws.Range("H113").Value = 9936.125
$ws.Range("I113").Value = 13796.6
$ws.Range("J113").Value = 3502
$ws.Range("K113").Value = 13796.6
$ws.Range("L113").Value = 3502
$ws.Range("M113").Value = -10542.6
$ws.Range("N113").Value = -10010

$ws.Range("H115").Value = 287099.44
$ws.Range("I115").Value = 287099.44
$ws.Range("K115").Value = 861298.3200000001
$ws.Range("M115").Value = -859731.3200000001

$ws.Range("H118").Value = 167928.17
$ws.Range("I118").Value = 167928.17
$ws.Range("K118").Value = 503784.51
$ws.Range("M118").Value = -502127.51

$ws.Range("H122").Value = 2019978.4
$ws.Range("I122").Value = 2116036.5
$ws.Range("J122").Value = 1429906.6
$ws.Range("K122").Value = 6348109.5
$ws.Range("L122").Value = 4289719.800000001
$ws.Range("M122").Value = -6345659.5
$ws.Range("N122").Value = -4294619.800000001

$ws.Range("H132").Value = 3498.4243
$ws.Range("I132").Value = 2086.348
$ws.Range("K132").Value = 6259.044
$ws.Range("M132").Value = -3729.044

$ws.Range("H138").Value = 2650.7253
$ws.Range("I138").Value = 1285.7858
$ws.Range("K138").Value = 3857.3574
$ws.Range("M138").Value = 1282.6426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1665.5555
$ws.Range("I45").Value = 1832.3334
$ws.Range("K45").Value = 1832.3334
$ws.Range("M45").Value = -1455.3334

$ws.Range("H61").Value = 3324.3901
$ws.Range("I61").Value = 1932.1923
$ws.Range("K61").Value = 1932.1923
$ws.Range("M61").Value = -1720.1923

$ws.Range("H122").Value = 2033.6342
$ws.Range("I122").Value = 1585
$ws.Range("J122").Value = 4212.7144
$ws.Range("K122").Value = 4755
$ws.Range("L122").Value = 12638.1432
$ws.Range("M122").Value = -2305
$ws.Range("N122").Value = -17538.1432

$ws.Range("H136").Value = 3324.3901
$ws.Range("I136").Value = 1932.1923
$ws.Range("K136").Value = 5796.5769
$ws.Range("M136").Value = -3246.5769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3829.2
$ws.Range("I62").Value = 3927.8572
$ws.Range("K62").Value = 3927.8572
$ws.Range("M62").Value = -3303.8572

$ws.Range("H64").Value = 36709
$ws.Range("J64").Value = 36709
$ws.Range("L64").Value = 36709
$ws.Range("N64").Value = -37205

$ws.Range("H65").Value = 3829.2
$ws.Range("I65").Value = 3927.8572
$ws.Range("K65").Value = 19639.286
$ws.Range("M65").Value = -16519.286

$ws.Range("H67").Value = 36709
$ws.Range("J67").Value = 36709
$ws.Range("L67").Value = 36709
$ws.Range("N67").Value = -38425

$ws.Range("H68").Value = 40445.875
$ws.Range("J68").Value = 40445.875
$ws.Range("L68").Value = 40445.875
$ws.Range("N68").Value = -41943.875

$ws.Range("H71").Value = 40445.875
$ws.Range("J71").Value = 40445.875
$ws.Range("L71").Value = 121337.625
$ws.Range("N71").Value = -128825.625

$ws.Range("H99").Value = 365341
$ws.Range("I99").Value = 780953.7
$ws.Range("J99").Value = 27655.688
$ws.Range("K99").Value = 780953.7
$ws.Range("L99").Value = 27655.688
$ws.Range("M99").Value = -779455.7
$ws.Range("N99").Value = -30651.688

$ws.Range("H126").Value = 365341
$ws.Range("I126").Value = 780953.7
$ws.Range("J126").Value = 27655.688
$ws.Range("K126").Value = 2342861.1
$ws.Range("L126").Value = 82967.064
$ws.Range("M126").Value = -2340391.1
$ws.Range("N126").Value = -87907.064

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 125071220
$ws.Range("J37").Value = 125071220
$ws.Range("L37").Value = 375213660
$ws.Range("N37").Value = -375213884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2318.075
$ws.Range("I97").Value = 1700.0667
$ws.Range("K97").Value = 1700.0667
$ws.Range("M97").Value = -1204.0667

$ws.Range("H132").Value = 22261.037
$ws.Range("I132").Value = 26362.791
$ws.Range("K132").Value = 79088.37300000001
$ws.Range("M132").Value = -76558.37300000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1368.2727
$ws.Range("I61").Value = 1285.6471
$ws.Range("J61").Value = 1649.2
$ws.Range("K61").Value = 1285.6471
$ws.Range("L61").Value = 1649.2
$ws.Range("M61").Value = -1083.6471
$ws.Range("N61").Value = -2053.2

$ws.Range("H68").Value = 3401.4814
$ws.Range("I68").Value = 3849.5264
$ws.Range("J68").Value = 2337.375
$ws.Range("K68").Value = 3849.5264
$ws.Range("L68").Value = 2337.375
$ws.Range("M68").Value = -3100.5264
$ws.Range("N68").Value = -3835.375

$ws.Range("H71").Value = 3401.4814
$ws.Range("I71").Value = 3849.5264
$ws.Range("J71").Value = 2337.375
$ws.Range("K71").Value = 19247.632
$ws.Range("L71").Value = 11686.875
$ws.Range("M71").Value = -15503.632
$ws.Range("N71").Value = -19174.875

$ws.Range("H113").Value = 1368.2727
$ws.Range("I113").Value = 1285.6471
$ws.Range("J113").Value = 1649.2
$ws.Range("K113").Value = 1285.6471
$ws.Range("L113").Value = 1649.2
$ws.Range("M113").Value = 884.3529000000001
$ws.Range("N113").Value = -5989.2

$ws.Range("H136").Value = 1152.8235
$ws.Range("I136").Value = 1152.8235
$ws.Range("K136").Value = 3458.4705
$ws.Range("M136").Value = -908.4704999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H132").Value = 2289.0952
$ws.Range("I132").Value = 2264.7
$ws.Range("K132").Value = 6794.099999999999
$ws.Range("M132").Value = -4264.099999999999
